# Automatically re-generate list and index
#
# The "Review date" column (B) in Sheet1 stores dates as literal text
# strings ("yyyy-mm-dd"). Every date whose day-of-month is "04" is
# bumped to day-of-month "06"; all other dates are left untouched.
#
# NumberFormat is forced to "@" (Text) before writing so Excel does not
# auto-convert the yyyy-mm-dd strings into date serials, then the
# range style is reset to "Normal" afterwards so the cells keep their
# original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$datesRange = $ws.Range("B2:B173")
$datesRange.NumberFormat = "@"

$ws.Range("B2").Value = "2011-03-06"
$ws.Range("B3").Value = "2011-06-06"
$ws.Range("B4").Value = "2011-07-06"
$ws.Range("B5").Value = "2011-12-06"
$ws.Range("B6").Value = "2014-09-06"
$ws.Range("B7").Value = "2014-10-06"
$ws.Range("B8").Value = "2014-11-06"
$ws.Range("B9").Value = "2014-12-06"
$ws.Range("B10").Value = "2014-12-06"
$ws.Range("B11").Value = "2015-03-06"
$ws.Range("B12").Value = "2016-01-06"
$ws.Range("B13").Value = "2016-04-06"
$ws.Range("B14").Value = "2016-08-06"
$ws.Range("B15").Value = "2017-01-06"
$ws.Range("B16").Value = "2017-08-06"
$ws.Range("B17").Value = "2017-09-06"
$ws.Range("B18").Value = "2017-09-06"
$ws.Range("B19").Value = "2017-09-06"
$ws.Range("B20").Value = "2017-10-06"
$ws.Range("B21").Value = "2017-10-06"
$ws.Range("B22").Value = "2017-11-06"
$ws.Range("B23").Value = "2018-05-06"
$ws.Range("B24").Value = "2018-07-06"
$ws.Range("B25").Value = "2019-05-06"
$ws.Range("B26").Value = "2019-05-06"
$ws.Range("B27").Value = "2019-06-06"
$ws.Range("B28").Value = "2019-06-06"
$ws.Range("B29").Value = "2019-06-06"
$ws.Range("B30").Value = "2019-08-06"
$ws.Range("B31").Value = "2019-09-06"
$ws.Range("B32").Value = "2020-01-06"
$ws.Range("B33").Value = "2020-03-06"
$ws.Range("B34").Value = "2020-03-06"
$ws.Range("B35").Value = "2020-03-06"
$ws.Range("B36").Value = "2020-03-06"
$ws.Range("B37").Value = "2020-03-06"
$ws.Range("B38").Value = "2020-03-06"
$ws.Range("B39").Value = "2020-04-06"
$ws.Range("B40").Value = "2020-05-06"
$ws.Range("B41").Value = "2020-07-06"
$ws.Range("B42").Value = "2020-08-06"
$ws.Range("B43").Value = "2020-11-06"
$ws.Range("B44").Value = "2021-02-06"
$ws.Range("B45").Value = "2021-05-06"
$ws.Range("B46").Value = "2021-05-06"
$ws.Range("B47").Value = "2021-05-06"
$ws.Range("B48").Value = "2021-06-06"
$ws.Range("B49").Value = "2021-06-06"
$ws.Range("B50").Value = "2021-06-06"
$ws.Range("B51").Value = "2021-06-06"
$ws.Range("B52").Value = "2021-06-06"
$ws.Range("B53").Value = "2021-06-06"
$ws.Range("B54").Value = "2021-06-06"
$ws.Range("B55").Value = "2021-06-06"
$ws.Range("B56").Value = "2021-06-06"
$ws.Range("B57").Value = "2021-09-06"
$ws.Range("B58").Value = "2021-12-06"
$ws.Range("B59").Value = "2021-12-06"
$ws.Range("B60").Value = "2022-01-06"
$ws.Range("B61").Value = "2022-01-06"
$ws.Range("B62").Value = "2022-03-06"
$ws.Range("B63").Value = "2022-03-06"
$ws.Range("B64").Value = "2022-04-06"
$ws.Range("B65").Value = "2022-05-06"
$ws.Range("B66").Value = "2022-05-06"
$ws.Range("B67").Value = "2022-06-06"
$ws.Range("B68").Value = "2022-06-06"
$ws.Range("B69").Value = "2022-06-06"
$ws.Range("B70").Value = "2022-06-06"
$ws.Range("B71").Value = "2022-07-06"
$ws.Range("B72").Value = "2022-07-06"
$ws.Range("B73").Value = "2022-07-06"
$ws.Range("B74").Value = "2022-08-06"
$ws.Range("B75").Value = "2022-09-06"
$ws.Range("B76").Value = "2022-10-06"
$ws.Range("B77").Value = "2022-10-06"
$ws.Range("B78").Value = "2022-11-06"
$ws.Range("B79").Value = "2022-11-06"
$ws.Range("B80").Value = "2022-12-06"
$ws.Range("B81").Value = "2022-12-06"
$ws.Range("B82").Value = "2022-12-06"
$ws.Range("B83").Value = "2023-01-06"
$ws.Range("B84").Value = "2023-01-06"
$ws.Range("B86").Value = "2023-02-06"
$ws.Range("B87").Value = "2023-02-06"
$ws.Range("B91").Value = "2023-04-06"
$ws.Range("B92").Value = "2023-04-06"
$ws.Range("B93").Value = "2023-04-06"
$ws.Range("B94").Value = "2023-04-06"
$ws.Range("B95").Value = "2023-05-06"
$ws.Range("B96").Value = "2023-05-06"
$ws.Range("B97").Value = "2023-05-06"
$ws.Range("B98").Value = "2023-05-06"
$ws.Range("B99").Value = "2023-05-06"
$ws.Range("B100").Value = "2023-05-06"
$ws.Range("B101").Value = "2023-05-06"
$ws.Range("B102").Value = "2023-05-06"
$ws.Range("B103").Value = "2023-05-06"
$ws.Range("B104").Value = "2023-05-06"
$ws.Range("B106").Value = "2023-06-06"
$ws.Range("B107").Value = "2023-06-06"
$ws.Range("B108").Value = "2023-06-06"
$ws.Range("B109").Value = "2023-06-06"
$ws.Range("B110").Value = "2023-06-06"
$ws.Range("B112").Value = "2023-07-06"
$ws.Range("B113").Value = "2023-07-06"
$ws.Range("B114").Value = "2023-07-06"
$ws.Range("B115").Value = "2023-07-06"
$ws.Range("B116").Value = "2023-07-06"
$ws.Range("B117").Value = "2023-08-06"
$ws.Range("B118").Value = "2023-08-06"
$ws.Range("B119").Value = "2023-08-06"
$ws.Range("B121").Value = "2023-09-06"
$ws.Range("B122").Value = "2023-09-06"
$ws.Range("B123").Value = "2023-10-06"
$ws.Range("B124").Value = "2023-10-06"
$ws.Range("B129").Value = "2024-01-06"
$ws.Range("B130").Value = "2024-02-06"
$ws.Range("B131").Value = "2024-02-06"
$ws.Range("B132").Value = "2024-03-06"
$ws.Range("B133").Value = "2024-03-06"
$ws.Range("B134").Value = "2024-03-06"
$ws.Range("B135").Value = "2024-03-06"
$ws.Range("B136").Value = "2024-03-06"
$ws.Range("B137").Value = "2024-04-06"
$ws.Range("B138").Value = "2024-04-06"
$ws.Range("B139").Value = "2024-05-06"
$ws.Range("B140").Value = "2024-06-06"
$ws.Range("B141").Value = "2024-07-06"
$ws.Range("B142").Value = "2024-07-06"
$ws.Range("B143").Value = "2024-07-06"
$ws.Range("B144").Value = "2024-07-06"
$ws.Range("B145").Value = "2024-07-06"
$ws.Range("B146").Value = "2024-08-06"
$ws.Range("B147").Value = "2024-08-06"
$ws.Range("B148").Value = "2024-10-06"
$ws.Range("B149").Value = "2024-10-06"
$ws.Range("B150").Value = "2025-01-06"
$ws.Range("B151").Value = "2025-01-06"
$ws.Range("B152").Value = "2025-02-06"
$ws.Range("B153").Value = "2025-02-06"
$ws.Range("B154").Value = "2025-02-06"
$ws.Range("B155").Value = "2025-03-06"
$ws.Range("B156").Value = "2025-03-06"
$ws.Range("B157").Value = "2025-04-06"
$ws.Range("B158").Value = "2025-05-06"
$ws.Range("B159").Value = "2025-05-06"
$ws.Range("B160").Value = "2025-05-06"
$ws.Range("B161").Value = "2025-05-06"
$ws.Range("B162").Value = "2025-06-06"
$ws.Range("B163").Value = "2025-07-06"
$ws.Range("B164").Value = "2025-07-06"
$ws.Range("B165").Value = "2025-08-06"
$ws.Range("B166").Value = "2025-10-06"
$ws.Range("B167").Value = "2025-10-06"
$ws.Range("B168").Value = "2026-01-06"
$ws.Range("B169").Value = "2026-03-06"
$ws.Range("B170").Value = "2026-03-06"
$ws.Range("B171").Value = "2026-03-06"
$ws.Range("B172").Value = "2027-01-06"
$ws.Range("B173").Value = "2027-02-06"

$datesRange.Style = "Normal"
